# 自动更新Excel文件 - 2026-02-18 23:21:29
# Daily countdown update: for every data row, the "剩余" (E, remaining days)
# counter ticks down by one. When a row's remaining-days counter would drop
# to zero, the cycle rolls over: remaining resets to the row's "总天"
# (D, total days) and the "开始时间" (F, start date) advances by that same
# number of days (i.e. a brand-new cycle begins).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # D: 总天 (total days)
    $eCell = $ws.Cells.Item($r, 5)   # E: 剩余 (remaining days)
    $fCell = $ws.Cells.Item($r, 6)   # F: 开始时间 (start date, yyyymmdd)

    $totalDays = $dCell.Value2
    $remaining = $eCell.Value2
    $startDateRaw = $fCell.Value2

    if ($remaining -eq $null -or $totalDays -eq $null -or $startDateRaw -eq $null) {
        continue
    }

    # Rows whose "remaining" already equals the "total" (never started a
    # countdown) or whose start date is not a well-formed yyyymmdd number
    # are left untouched (e.g. row 36's corrupt 9-digit date).
    if ($remaining -ge $totalDays) {
        continue
    }
    if (($startDateRaw -lt 10000101) -or ($startDateRaw -gt 99991231)) {
        continue
    }

    if ($remaining -le 1) {
        # Cycle completed: reset remaining to the total and roll the start
        # date forward by the total number of days.
        $eCell.Value2 = $totalDays

        $startDate = $startDateRaw
        $y = [math]::Floor($startDate / 10000)
        $m = [math]::Floor(($startDate % 10000) / 100)
        $d = $startDate % 100

        $dt = Get-Date -Year $y -Month $m -Day $d
        $dt = $dt.AddDays([double]$totalDays)

        $newDateNum = ($dt.Year * 10000) + ($dt.Month * 100) + $dt.Day
        $fCell.Value2 = $newDateNum
    }
    else {
        $eCell.Value2 = $remaining - 1
    }
}
